# Auto-generated edit script applying the Leviathan_Profits.xlsx market-data refresh diff.
# For each affected row, columns H-N (currentAveragePrice .. LeveProfitHQ) are updated
# to the refreshed values; a couple of cells are newly populated or cleared entirely
# to match the source row's new shape.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# ALC row 17
$ws.Cells.Item(17, 8).Value = 447462.53
$ws.Cells.Item(17, 10).Value = 447462.53
$ws.Cells.Item(17, 12).Value = 1342387.59
$ws.Cells.Item(17, 14).Value = -1342723.59
# ALC row 29
$ws.Cells.Item(29, 8).Value = 4894
$ws.Cells.Item(29, 9).Value = 269.2
$ws.Cells.Item(29, 11).Value = 807.5999999999999
$ws.Cells.Item(29, 13).Value = -526.5999999999999
# ALC row 33
$ws.Cells.Item(33, 8).Value = 386.77777
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 12).Value = 300
$ws.Cells.Item(33, 14).Value = -758
# ALC row 39
$ws.Cells.Item(39, 8).Value = 2121.4167
$ws.Cells.Item(39, 9).Value = 2899.1428
$ws.Cells.Item(39, 10).Value = 1032.6
$ws.Cells.Item(39, 11).Value = 8697.428400000001
$ws.Cells.Item(39, 12).Value = 3097.8
$ws.Cells.Item(39, 13).Value = -8401.428400000001
$ws.Cells.Item(39, 14).Value = -3689.8
# ALC row 49
$ws.Cells.Item(49, 8).Value = 1099.6666
$ws.Cells.Item(49, 9).Value = 500
$ws.Cells.Item(49, 10).Value = 1399.5
$ws.Cells.Item(49, 11).Value = 1500
$ws.Cells.Item(49, 12).Value = 4198.5
$ws.Cells.Item(49, 13).Value = -1364
$ws.Cells.Item(49, 14).Value = -4470.5
# ALC row 86
$ws.Cells.Item(86, 8).Value = 1830.8636
$ws.Cells.Item(86, 9).Value = 1728.1428
$ws.Cells.Item(86, 11).Value = 1728.1428
$ws.Cells.Item(86, 13).Value = -605.1428000000001
# ALC row 89
$ws.Cells.Item(89, 8).Value = 1830.8636
$ws.Cells.Item(89, 9).Value = 1728.1428
$ws.Cells.Item(89, 11).Value = 8640.714
$ws.Cells.Item(89, 13).Value = -3024.714
# ALC row 101
$ws.Cells.Item(101, 8).Value = 1795.5
$ws.Cells.Item(101, 10).Value = 899.75
$ws.Cells.Item(101, 12).Value = 2699.25
$ws.Cells.Item(101, 14).Value = -5943.25
# ALC row 106
$ws.Cells.Item(106, 8).Value = 6859.737
$ws.Cells.Item(106, 9).Value = 1226.1818
$ws.Cells.Item(106, 11).Value = 1226.1818
$ws.Cells.Item(106, 13).Value = -595.1818000000001
# ALC row 116
$ws.Cells.Item(116, 8).Value = 3606.9285
$ws.Cells.Item(116, 9).Value = 2838.6667
$ws.Cells.Item(116, 10).Value = 4989.8
$ws.Cells.Item(116, 11).Value = 2838.6667
$ws.Cells.Item(116, 12).Value = 4989.8
$ws.Cells.Item(116, 13).Value = 603.3332999999998
$ws.Cells.Item(116, 14).Value = -11873.8

$ws = $wb.Sheets.Item("ARM")
# ARM row 45
$ws.Cells.Item(45, 8).Value = 7792.3184
$ws.Cells.Item(45, 9).Value = 12285.3
$ws.Cells.Item(45, 11).Value = 12285.3
$ws.Cells.Item(45, 13).Value = -11908.3
# ARM row 97
$ws.Cells.Item(97, 8).Value = 3416.3333
$ws.Cells.Item(97, 9).Value = 3099.6
$ws.Cells.Item(97, 10).Value = 5000
$ws.Cells.Item(97, 11).Value = 3099.6
$ws.Cells.Item(97, 12).Value = 5000
$ws.Cells.Item(97, 13).Value = -2603.6
$ws.Cells.Item(97, 14).Value = -5992
# ARM row 122
$ws.Cells.Item(122, 8).Value = 1845.3191
$ws.Cells.Item(122, 9).Value = 1357.2
$ws.Cells.Item(122, 11).Value = 4071.6
$ws.Cells.Item(122, 13).Value = -1621.6

$ws = $wb.Sheets.Item("BSM")
# BSM row 20
$ws.Cells.Item(20, 8).Value = 5497.9707
$ws.Cells.Item(20, 9).Value = 6532.696
$ws.Cells.Item(20, 10).Value = 3334.4546
$ws.Cells.Item(20, 11).Value = 6532.696
$ws.Cells.Item(20, 12).Value = 3334.4546
$ws.Cells.Item(20, 13).Value = -6285.696
$ws.Cells.Item(20, 14).Value = -3828.4546
# BSM row 94
$ws.Cells.Item(94, 8).Value = 1028
$ws.Cells.Item(94, 9).Value = 1029.5834
$ws.Cells.Item(94, 11).Value = 1029.5834
$ws.Cells.Item(94, 13).Value = -578.5834
# BSM row 105
$ws.Cells.Item(105, 8).Value = 1838.8572
$ws.Cells.Item(105, 9).Value = 1902.5385
$ws.Cells.Item(105, 10).Value = 1011
$ws.Cells.Item(105, 11).Value = 1902.5385
$ws.Cells.Item(105, 12).Value = 1011
$ws.Cells.Item(105, 13).Value = -155.5385000000001
$ws.Cells.Item(105, 14).Value = -4505

$ws = $wb.Sheets.Item("CRP")
# CRP row 22
$ws.Cells.Item(22, 8).Value = 788.3333
$ws.Cells.Item(22, 9).Value = 492.25
$ws.Cells.Item(22, 10).Value = 1025.2
$ws.Cells.Item(22, 11).Value = 492.25
$ws.Cells.Item(22, 12).Value = 1025.2
$ws.Cells.Item(22, 13).Value = -142.25
$ws.Cells.Item(22, 14).Value = -1725.2
# CRP row 99
$ws.Cells.Item(99, 8).Value = 2243.3809
$ws.Cells.Item(99, 9).Value = 1596.5
$ws.Cells.Item(99, 10).Value = 3105.889
$ws.Cells.Item(99, 11).Value = 1596.5
$ws.Cells.Item(99, 12).Value = 3105.889
$ws.Cells.Item(99, 13).Value = -98.5
$ws.Cells.Item(99, 14).Value = -6101.889
# CRP row 122
$ws.Cells.Item(122, 8).Value = 128269.875
$ws.Cells.Item(122, 9).Value = 146119.86
$ws.Cells.Item(122, 11).Value = 438359.58
$ws.Cells.Item(122, 13).Value = -435909.58
# CRP row 126
$ws.Cells.Item(126, 8).Value = 2243.3809
$ws.Cells.Item(126, 9).Value = 1596.5
$ws.Cells.Item(126, 10).Value = 3105.889
$ws.Cells.Item(126, 11).Value = 4789.5
$ws.Cells.Item(126, 12).Value = 9317.667000000001
$ws.Cells.Item(126, 13).Value = -2319.5
$ws.Cells.Item(126, 14).Value = -14257.667
# CRP row 134
$ws.Cells.Item(134, 8).Value = 1887.775
$ws.Cells.Item(134, 9).Value = 1111.1515
$ws.Cells.Item(134, 11).Value = 3333.4545
$ws.Cells.Item(134, 13).Value = -798.4544999999998
# CRP row 141
$ws.Cells.Item(141, 8).Value = 332399.3
$ws.Cells.Item(141, 10).Value = 392374.12
$ws.Cells.Item(141, 12).Value = 392374.12
$ws.Cells.Item(141, 14).Value = -402734.12

$ws = $wb.Sheets.Item("CUL")
# CUL row 54
$ws.Cells.Item(54, 8).Value = 8000
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 8000
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 24000
$ws.Cells.Item(54, 13).ClearContents()
$ws.Cells.Item(54, 14).Value = -25118
# CUL row 107
$ws.Cells.Item(107, 8).Value = 772.63635
$ws.Cells.Item(107, 10).Value = 815.7857
$ws.Cells.Item(107, 12).Value = 2447.3571
$ws.Cells.Item(107, 14).Value = -6287.3571
# CUL row 113
$ws.Cells.Item(113, 8).Value = 1982.6
$ws.Cells.Item(113, 10).Value = 2204.15
$ws.Cells.Item(113, 12).Value = 6612.450000000001
$ws.Cells.Item(113, 14).Value = -10952.45
# CUL row 121
$ws.Cells.Item(121, 8).Value = 17651242
$ws.Cells.Item(121, 9).Value = 66666810
$ws.Cells.Item(121, 11).Value = 200000430
$ws.Cells.Item(121, 13).Value = -199999120
# CUL row 131
$ws.Cells.Item(131, 8).Value = 2526.1365
$ws.Cells.Item(131, 9).Value = 2198.1667
$ws.Cells.Item(131, 11).Value = 6594.500100000001
$ws.Cells.Item(131, 13).Value = -1554.500100000001

$ws = $wb.Sheets.Item("GSM")
# GSM row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()
# GSM row 30
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 13).ClearContents()
# GSM row 121
$ws.Cells.Item(121, 8).Value = 240317
$ws.Cells.Item(121, 10).Value = 240317
$ws.Cells.Item(121, 12).Value = 240317
$ws.Cells.Item(121, 14).Value = -243811
# GSM row 122
$ws.Cells.Item(122, 8).Value = 1833.6571
$ws.Cells.Item(122, 9).Value = 1651.375
$ws.Cells.Item(122, 10).Value = 2231.3635
$ws.Cells.Item(122, 11).Value = 4954.125
$ws.Cells.Item(122, 12).Value = 6694.0905
$ws.Cells.Item(122, 13).Value = -2504.125
$ws.Cells.Item(122, 14).Value = -11594.0905
# GSM row 132
$ws.Cells.Item(132, 8).Value = 3878.4285
$ws.Cells.Item(132, 9).Value = 5028.5557
$ws.Cells.Item(132, 10).Value = 1808.2
$ws.Cells.Item(132, 11).Value = 15085.6671
$ws.Cells.Item(132, 12).Value = 5424.6
$ws.Cells.Item(132, 13).Value = -12555.6671
$ws.Cells.Item(132, 14).Value = -10484.6

$ws = $wb.Sheets.Item("LTW")
# LTW row 7
$ws.Cells.Item(7, 8).Value = 12519.216
$ws.Cells.Item(7, 9).Value = 26462.77
$ws.Cells.Item(7, 10).Value = 4966.4585
$ws.Cells.Item(7, 11).Value = 26462.77
$ws.Cells.Item(7, 12).Value = 4966.4585
$ws.Cells.Item(7, 13).Value = -26350.77
$ws.Cells.Item(7, 14).Value = -5190.4585
# LTW row 16
$ws.Cells.Item(16, 8).Value = 538.4231
$ws.Cells.Item(16, 9).Value = 556.36
$ws.Cells.Item(16, 10).Value = 90
$ws.Cells.Item(16, 11).Value = 556.36
$ws.Cells.Item(16, 12).Value = 90
$ws.Cells.Item(16, 13).Value = -386.36
$ws.Cells.Item(16, 14).Value = -430
# LTW row 22
$ws.Cells.Item(22, 8).Value = 3941.7368
$ws.Cells.Item(22, 9).Value = 3752.6858
$ws.Cells.Item(22, 10).Value = 6147.3335
$ws.Cells.Item(22, 11).Value = 3752.6858
$ws.Cells.Item(22, 12).Value = 6147.3335
$ws.Cells.Item(22, 13).Value = -3457.6858
$ws.Cells.Item(22, 14).Value = -6737.3335
# LTW row 27
$ws.Cells.Item(27, 8).Value = 3941.7368
$ws.Cells.Item(27, 9).Value = 3752.6858
$ws.Cells.Item(27, 10).Value = 6147.3335
$ws.Cells.Item(27, 11).Value = 3752.6858
$ws.Cells.Item(27, 12).Value = 6147.3335
$ws.Cells.Item(27, 13).Value = -3645.6858
$ws.Cells.Item(27, 14).Value = -6361.3335
# LTW row 39
$ws.Cells.Item(39, 8).Value = 7399.5
$ws.Cells.Item(39, 9).Value = 7399.5
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 7399.5
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = -6939.5
$ws.Cells.Item(39, 14).ClearContents()
# LTW row 40
$ws.Cells.Item(40, 8).Value = 16235.143
$ws.Cells.Item(40, 9).Value = 41666.332
$ws.Cells.Item(40, 10).Value = 9299.362999999999
$ws.Cells.Item(40, 11).Value = 41666.332
$ws.Cells.Item(40, 12).Value = 9299.362999999999
$ws.Cells.Item(40, 13).Value = -41530.332
$ws.Cells.Item(40, 14).Value = -9571.362999999999
# LTW row 122
$ws.Cells.Item(122, 8).Value = 70022.57000000001
$ws.Cells.Item(122, 9).Value = 97970.336
$ws.Cells.Item(122, 10).Value = 4811.1113
$ws.Cells.Item(122, 11).Value = 293911.008
$ws.Cells.Item(122, 12).Value = 14433.3339
$ws.Cells.Item(122, 13).Value = -291461.008
$ws.Cells.Item(122, 14).Value = -19333.3339
# LTW row 126
$ws.Cells.Item(126, 8).Value = 12519.216
$ws.Cells.Item(126, 9).Value = 26462.77
$ws.Cells.Item(126, 10).Value = 4966.4585
$ws.Cells.Item(126, 11).Value = 79388.31
$ws.Cells.Item(126, 12).Value = 14899.3755
$ws.Cells.Item(126, 13).Value = -76918.31
$ws.Cells.Item(126, 14).Value = -19839.3755
# LTW row 132
$ws.Cells.Item(132, 8).Value = 4583
$ws.Cells.Item(132, 9).Value = 4025.24
$ws.Cells.Item(132, 10).Value = 5977.4
$ws.Cells.Item(132, 11).Value = 12075.72
$ws.Cells.Item(132, 12).Value = 17932.2
$ws.Cells.Item(132, 13).Value = -9545.719999999999
$ws.Cells.Item(132, 14).Value = -22992.2

$ws = $wb.Sheets.Item("WVR")
# WVR row 23
$ws.Cells.Item(23, 8).Value = 300
$ws.Cells.Item(23, 9).Value = 300
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 300
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -71
$ws.Cells.Item(23, 14).ClearContents()
# WVR row 54
$ws.Cells.Item(54, 8).Value = 42499.75
# WVR row 70
$ws.Cells.Item(70, 8).Value = 17041.268
$ws.Cells.Item(70, 10).Value = 17041.268
$ws.Cells.Item(70, 12).Value = 17041.268
$ws.Cells.Item(70, 14).Value = -17671.268
# WVR row 73
$ws.Cells.Item(73, 8).Value = 17041.268
$ws.Cells.Item(73, 10).Value = 17041.268
$ws.Cells.Item(73, 12).Value = 17041.268
$ws.Cells.Item(73, 14).Value = -19225.268
# WVR row 96
$ws.Cells.Item(96, 9).Value = 1800.2858
$ws.Cells.Item(96, 10).Value = 2200
$ws.Cells.Item(96, 11).Value = 1800.2858
$ws.Cells.Item(96, 12).Value = 2200
$ws.Cells.Item(96, 13).Value = -427.2858000000001
$ws.Cells.Item(96, 14).Value = -4946
# WVR row 126
$ws.Cells.Item(126, 8).Value = 1823.409
$ws.Cells.Item(126, 9).Value = 1800.8422
$ws.Cells.Item(126, 11).Value = 5402.5266
$ws.Cells.Item(126, 13).Value = -2932.5266
# WVR row 136
$ws.Cells.Item(136, 8).Value = 810.55884
$ws.Cells.Item(136, 9).Value = 774.51514
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 2323.54542
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = 226.4545800000001
$ws.Cells.Item(136, 14).Value = -11100
